# Scheduled market-data refresh: updates currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) for specific Leve rows across the per-job worksheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR). Values sourced from the latest Universalis market snapshot.

$wb = $excel.ActiveWorkbook

# ===================== Sheet: ALC =====================
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 170
$ws.Range("I12").Value = 170
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 170
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0
$ws.Range("N12").ClearContents()

# Row 17
$ws.Range("H17").Value = 2687.4285
$ws.Range("I17").Value = 2147.5
$ws.Range("J17").Value = 2777.4167
$ws.Range("K17").Value = 6442.5
$ws.Range("L17").Value = 8332.250100000001
$ws.Range("M17").Value = -6274.5
$ws.Range("N17").Value = -8668.250100000001

# Row 32
$ws.Range("H32").Value = 3290.4285
$ws.Range("J32").Value = 3406.6
$ws.Range("L32").Value = 3406.6
$ws.Range("N32").Value = -4058.6

# Row 40
$ws.Range("H40").Value = 1858.7273
$ws.Range("J40").Value = 1600
$ws.Range("L40").Value = 1600
$ws.Range("N40").Value = -1950

# Row 43
$ws.Range("H43").Value = 2199.8
$ws.Range("I43").Value = 1500
$ws.Range("J43").Value = 2666.3333
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 2666.3333
$ws.Range("M43").Value = -1431
$ws.Range("N43").Value = -2804.3333

# Row 62
$ws.Range("H62").Value = 5875
$ws.Range("I62").Value = 3750
$ws.Range("K62").Value = 3750
$ws.Range("M62").Value = -3126

# Row 65
$ws.Range("H65").Value = 5875
$ws.Range("I65").Value = 3750
$ws.Range("K65").Value = 18750
$ws.Range("M65").Value = -15630

# Row 75
$ws.Range("H75").Value = 33590
$ws.Range("J75").Value = 33590
$ws.Range("L75").Value = 33590
$ws.Range("N75").Value = -35462

# Row 78
$ws.Range("H78").Value = 33590
$ws.Range("J78").Value = 33590
$ws.Range("L78").Value = 100770
$ws.Range("N78").Value = -110130

# Row 100
$ws.Range("H100").Value = 2455.889
$ws.Range("J100").Value = 5000
$ws.Range("L100").Value = 5000
$ws.Range("N100").Value = -6082

# Row 111
$ws.Range("H111").Value = 1514.5
$ws.Range("I111").Value = 529
$ws.Range("K111").Value = 1587
$ws.Range("M111").Value = 1480

# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# Row 125
$ws.Range("H125").Value = 214290700
$ws.Range("J125").Value = 166667020
$ws.Range("L125").Value = 1500003180
$ws.Range("N125").Value = -1500008100

# Row 132
$ws.Range("H132").Value = 1225.6666
$ws.Range("I132").Value = 1270.8
$ws.Range("K132").Value = 3812.4
$ws.Range("M132").Value = -1282.4

# Row 138
$ws.Range("H138").Value = 2602.442
$ws.Range("I138").Value = 2774.5908
$ws.Range("J138").Value = 2422.0952
$ws.Range("K138").Value = 8323.7724
$ws.Range("L138").Value = 7266.285600000001
$ws.Range("M138").Value = -3183.7724
$ws.Range("N138").Value = -17546.2856


# ===================== Sheet: ARM =====================
$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 4279.75
$ws.Range("I61").Value = 835
$ws.Range("J61").Value = 7724.5
$ws.Range("K61").Value = 835
$ws.Range("L61").Value = 7724.5
$ws.Range("M61").Value = -623
$ws.Range("N61").Value = -8148.5

# Row 136
$ws.Range("H136").Value = 4279.75
$ws.Range("I136").Value = 835
$ws.Range("J136").Value = 7724.5
$ws.Range("K136").Value = 2505
$ws.Range("L136").Value = 23173.5
$ws.Range("M136").Value = 45
$ws.Range("N136").Value = -28273.5


# ===================== Sheet: BSM =====================
$ws = $wb.Worksheets.Item("BSM")

# Row 99
$ws.Range("H99").Value = 1498
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1498
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1498
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4494

# Row 105
$ws.Range("H105").Value = 2849.7
$ws.Range("I105").Value = 2785.4285
$ws.Range("K105").Value = 2785.4285
$ws.Range("M105").Value = -1038.4285

# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# Row 134
$ws.Range("H134").Value = 893.75
$ws.Range("I134").Value = 893
$ws.Range("K134").Value = 2679
$ws.Range("M134").Value = -144


# ===================== Sheet: CRP =====================
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 4077.25
$ws.Range("I16").Value = 4437
$ws.Range("K16").Value = 4437
$ws.Range("M16").Value = -4150

# Row 22
$ws.Range("H22").Value = 355.77777
$ws.Range("I22").Value = 366.33334
$ws.Range("J22").Value = 334.66666
$ws.Range("K22").Value = 366.33334
$ws.Range("L22").Value = 334.66666
$ws.Range("M22").Value = -16.33334000000002
$ws.Range("N22").Value = -1034.66666

# Row 31
$ws.Range("H31").Value = 2090.182
$ws.Range("I31").Value = 1638
$ws.Range("J31").Value = 2348.5715
$ws.Range("K31").Value = 1638
$ws.Range("L31").Value = 2348.5715
$ws.Range("M31").Value = -1343
$ws.Range("N31").Value = -2938.5715

# Row 34
$ws.Range("H34").Value = 2090.182
$ws.Range("I34").Value = 1638
$ws.Range("J34").Value = 2348.5715
$ws.Range("K34").Value = 1638
$ws.Range("L34").Value = 2348.5715
$ws.Range("M34").Value = -1436
$ws.Range("N34").Value = -2752.5715

# Row 113
$ws.Range("H113").Value = 4077.25
$ws.Range("I113").Value = 4437
$ws.Range("K113").Value = 4437
$ws.Range("M113").Value = -2267

# Row 133
$ws.Range("H133").Value = 45246

# Row 134
$ws.Range("H134").Value = 1217.5
$ws.Range("I134").Value = 1217.5
$ws.Range("K134").Value = 3652.5
$ws.Range("M134").Value = -1117.5


# ===================== Sheet: CUL =====================
$ws = $wb.Worksheets.Item("CUL")

# Row 2
$ws.Range("H2").Value = 138241.81
$ws.Range("I2").Value = 122258.11
$ws.Range("J2").Value = 158792.28
$ws.Range("K2").Value = 733548.66
$ws.Range("L2").Value = 952753.6799999999
$ws.Range("M2").Value = -733435.66
$ws.Range("N2").Value = -952979.6799999999

# Row 38
$ws.Range("H38").Value = 68.5
$ws.Range("I38").Value = 68.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 205.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 141.5
$ws.Range("N38").ClearContents()

# Row 41
$ws.Range("H41").Value = 3970.4707
$ws.Range("I41").Value = 880
$ws.Range("J41").Value = 4382.533
$ws.Range("K41").Value = 2640
$ws.Range("L41").Value = 13147.599
$ws.Range("M41").Value = -2302
$ws.Range("N41").Value = -13823.599

# Row 60
$ws.Range("H60").Value = 1096.5
$ws.Range("I60").Value = 299
$ws.Range("K60").Value = 897
$ws.Range("M60").Value = -646

# Row 122
$ws.Range("H122").Value = 892
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 892
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8028
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -12928


# ===================== Sheet: GSM =====================
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 7011
$ws.Range("I70").Value = 6900.9375
$ws.Range("J70").Value = 7171.091
$ws.Range("K70").Value = 6900.9375
$ws.Range("L70").Value = 7171.091
$ws.Range("M70").Value = -6630.9375
$ws.Range("N70").Value = -7711.091

# Row 73
$ws.Range("H73").Value = 7011
$ws.Range("I73").Value = 6900.9375
$ws.Range("J73").Value = 7171.091
$ws.Range("K73").Value = 6900.9375
$ws.Range("L73").Value = 7171.091
$ws.Range("M73").Value = -5964.9375
$ws.Range("N73").Value = -9043.091

# Row 122
$ws.Range("H122").Value = 2633.3333
$ws.Range("I122").Value = 2950
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8850
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6400
$ws.Range("N122").Value = -10900


# ===================== Sheet: LTW =====================
$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Range("H46").Value = 2751
$ws.Range("J46").Value = 2902
$ws.Range("L46").Value = 2902
$ws.Range("N46").Value = -3278

# Row 61
$ws.Range("H61").Value = 1640.6154
$ws.Range("J61").Value = 2154.125
$ws.Range("L61").Value = 2154.125
$ws.Range("N61").Value = -2558.125

# Row 82
$ws.Range("H82").Value = 3798.5
$ws.Range("I82").Value = 4731.3335
$ws.Range("J82").Value = 1000
$ws.Range("K82").Value = 4731.3335
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = -4370.3335
$ws.Range("N82").Value = -1722

# Row 85
$ws.Range("H85").Value = 3798.5
$ws.Range("I85").Value = 4731.3335
$ws.Range("J85").Value = 1000
$ws.Range("K85").Value = 4731.3335
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = -3483.3335
$ws.Range("N85").Value = -3496

# Row 113
$ws.Range("H113").Value = 1640.6154
$ws.Range("J113").Value = 2154.125
$ws.Range("L113").Value = 2154.125
$ws.Range("N113").Value = -6494.125

# Row 122
$ws.Range("H122").Value = 4747.8604
$ws.Range("I122").Value = 3785.3103
$ws.Range("K122").Value = 11355.9309
$ws.Range("M122").Value = -8905.930899999999

# Row 134
$ws.Range("H134").Value = 120000
$ws.Range("J134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -130140


# ===================== Sheet: WVR =====================
$ws = $wb.Worksheets.Item("WVR")

# Row 13
$ws.Range("H13").Value = 175
$ws.Range("I13").Value = 50
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 50
$ws.Range("L13").Value = 300
$ws.Range("M13").Value = 90
$ws.Range("N13").Value = -580

# Row 113
$ws.Range("H113").Value = 1946.2
$ws.Range("I113").Value = 1947
$ws.Range("J113").Value = 1943
$ws.Range("K113").Value = 5841
$ws.Range("L113").Value = 5829
$ws.Range("M113").Value = -3671
$ws.Range("N113").Value = -10169

# Row 126
$ws.Range("H126").Value = 2027.1666
$ws.Range("I126").Value = 1965.0714
$ws.Range("K126").Value = 5895.2142
$ws.Range("M126").Value = -3425.2142

# Row 132
$ws.Range("H132").Value = 4052.1765
$ws.Range("I132").Value = 4092.625
$ws.Range("J132").Value = 3405
$ws.Range("K132").Value = 12277.875
$ws.Range("L132").Value = 10215
$ws.Range("M132").Value = -9747.875
$ws.Range("N132").Value = -15275

# Row 136
$ws.Range("H136").Value = 3090.111
$ws.Range("I136").Value = 3083.524
$ws.Range("J136").Value = 3113.1667
$ws.Range("K136").Value = 9250.572
$ws.Range("L136").Value = 9339.500100000001
$ws.Range("M136").Value = -6700.572
$ws.Range("N136").Value = -14439.5001

